$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.291.93"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.694.79"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.67"
$ws.Range("E5").Value = "  +4.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.87"
$ws.Range("E6").Value = "  +13.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.719"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.40"
$ws.Range("E10").Value = "  +15.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.158"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000285"
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.35"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.276.09"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.690.01"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.29"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("E18").Value = "  +2.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.73"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.218.62"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "406.31"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "89.65"
$ws.Range("E23").Value = "  +2.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.63"
$ws.Range("E24").Value = "  +7.38%  "
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.91"
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.75"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.56"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.42"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.61"
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "47.50"
$ws.Range("E32").Value = "  +10.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.64"
$ws.Range("E33").Value = "  +2.39%  "
$ws.Range("E34").Value = "  +4.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "629.61"
$ws.Range("E35").Value = "  +6.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "67.28"
$ws.Range("E36").Value = "  +4.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0821"
$ws.Range("E37").Value = "  -7.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.409"
$ws.Range("E38").Value = "  +3.86%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.99"
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0441"
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.62"
$ws.Range("E44").Value = "  -4.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.877.18"
$ws.Range("E45").Value = "  +4.79%  "
$ws.Range("E46").Value = "  +4.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.18"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "145.41"
$ws.Range("E48").Value = "  +3.44%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.60"
$ws.Range("E50").Value = "  -6.61%  "
$ws.Range("E51").Value = "  -2.59%  "
